$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Delete rows 3 and 4 (005146441/JOSE and 004550605/REJANE) - delete row 4 first
# so the row number for the earlier deletion stays valid.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(3).Delete()

# After the two deletions above, former row 6 (005002457/ROSANGELA) is now row 4.
# Update its Saldo value from 20231.7 to 17000.
$ws.Cells.Item(4, 3).Value = 17000

# Former rows 7 (004748803/DORIVAL) and 8 (004202332/TATIANA) are now rows 5 and 6.
# Delete row 6 first, then row 5.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Former row 9 (005366255/RAPHAELA) is now row 5. Insert a new row after it
# for 004231509/THEOMAR/953.09. The account number has a leading zero, so
# force text formatting while entering it (otherwise Excel would store it
# as the number 4231509), then strip the number format back off so the
# cell ends up a plain text value with no special formatting applied.
$ws.Rows.Item(6).Insert()
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "004231509"
$ws.Cells.Item(6, 1).ClearFormats()
$ws.Cells.Item(6, 2).Value = "THEOMAR"
$ws.Cells.Item(6, 3).Value = 953.09

# Former row 15 (004474776/GILSON) is now row 12. Insert a new row after it
# for 005003629/ANDRE/381.23 (same leading-zero text handling as above).
$ws.Rows.Item(13).Insert()
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "005003629"
$ws.Cells.Item(13, 1).ClearFormats()
$ws.Cells.Item(13, 2).Value = "ANDRE"
$ws.Cells.Item(13, 3).Value = 381.23
